$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update each data row (rows 2-17) with new NATMI-computed values
# following Dr Hou advice: Ligand-expressing / Receptor-expressing cell
# counts changed from 1 to 3, and downstream expression/specificity
# values were recalculated accordingly.

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 122.253015
$ws.Range("H2").Value = 366.759045
$ws.Range("I2").Value = 0.1988639364328829
$ws.Range("J2").Value = 0.1988639364328829
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 10.50918
$ws.Range("N2").Value = 31.52754
$ws.Range("O2").Value = 0.1224541750083835
$ws.Range("P2").Value = 0.1224541750083835
$ws.Range("Q2").Value = 1284.7789401777
$ws.Range("R2").Value = 11563.0104615993
$ws.Range("S2").Value = 0.02435171927480829
$ws.Range("T2").Value = 0.02435171927480829

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 122.253015
$ws.Range("H3").Value = 366.759045
$ws.Range("I3").Value = 0.1988639364328829
$ws.Range("J3").Value = 0.1988639364328829
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.176377666666667
$ws.Range("N3").Value = 15.529133
$ws.Range("O3").Value = 0.06031574839364136
$ws.Range("P3").Value = 0.06031574839364136
$ws.Range("Q3").Value = 632.8277765286651
$ws.Range("R3").Value = 5695.449988757986
$ws.Range("S3").Value = 0.01199462715445485
$ws.Range("T3").Value = 0.01199462715445485

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 122.253015
$ws.Range("H4").Value = 366.759045
$ws.Range("I4").Value = 0.1988639364328829
$ws.Range("J4").Value = 0.1988639364328829
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 66.063113
$ws.Range("N4").Value = 198.189339
$ws.Range("O4").Value = 0.7697749968028538
$ws.Range("P4").Value = 0.769774996802854
$ws.Range("Q4").Value = 8076.414744535696
$ws.Range("R4").Value = 72687.73270082126
$ws.Range("S4").Value = 0.1530804860318253
$ws.Range("T4").Value = 0.1530804860318254

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 122.253015
$ws.Range("H5").Value = 366.759045
$ws.Range("I5").Value = 0.1988639364328829
$ws.Range("J5").Value = 0.1988639364328829
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.072658
$ws.Range("N5").Value = 12.217974
$ws.Range("O5").Value = 0.04745507979512132
$ws.Range("P5").Value = 0.04745507979512132
$ws.Range("Q5").Value = 497.89471956387
$ws.Range("R5").Value = 4481.052476074829
$ws.Range("S5").Value = 0.009437103971794392
$ws.Range("T5").Value = 0.00943710397179439

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 132.5447616666667
$ws.Range("H6").Value = 397.634285
$ws.Range("I6").Value = 0.2156050961899926
$ws.Range("J6").Value = 0.2156050961899926
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 10.50918
$ws.Range("N6").Value = 31.52754
$ws.Range("O6").Value = 0.1224541750083835
$ws.Range("P6").Value = 0.1224541750083835
$ws.Range("Q6").Value = 1392.9367584121
$ws.Range("R6").Value = 12536.4308257089
$ws.Range("S6").Value = 0.02640174418154871
$ws.Range("T6").Value = 0.02640174418154871

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 132.5447616666667
$ws.Range("H7").Value = 397.634285
$ws.Range("I7").Value = 0.2156050961899926
$ws.Range("J7").Value = 0.2156050961899926
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 5.176377666666667
$ws.Range("N7").Value = 15.529133
$ws.Range("O7").Value = 0.06031574839364136
$ws.Range("P7").Value = 0.06031574839364136
$ws.Range("Q7").Value = 686.1017441249895
$ws.Range("R7").Value = 6174.915697124905
$ws.Range("S7").Value = 0.01300438273418244
$ws.Range("T7").Value = 0.01300438273418244

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 132.5447616666667
$ws.Range("H8").Value = 397.634285
$ws.Range("I8").Value = 0.2156050961899926
$ws.Range("J8").Value = 0.2156050961899926
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 66.063113
$ws.Range("N8").Value = 198.189339
$ws.Range("O8").Value = 0.7697749968028538
$ws.Range("P8").Value = 0.769774996802854
$ws.Range("Q8").Value = 8756.319567543069
$ws.Range("R8").Value = 78806.87610788761
$ws.Range("S8").Value = 0.1659674122303305
$ws.Range("T8").Value = 0.1659674122303306

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 132.5447616666667
$ws.Range("H9").Value = 397.634285
$ws.Range("I9").Value = 0.2156050961899926
$ws.Range("J9").Value = 0.2156050961899926
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.072658
$ws.Range("N9").Value = 12.217974
$ws.Range("O9").Value = 0.04745507979512132
$ws.Range("P9").Value = 0.04745507979512132
$ws.Range("Q9").Value = 539.8094839598433
$ws.Range("R9").Value = 4858.285355638589
$ws.Range("S9").Value = 0.01023155704393091
$ws.Range("T9").Value = 0.01023155704393091

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 320.0894206666666
$ws.Range("H10").Value = 960.2682619999999
$ws.Range("I10").Value = 0.5206762565675317
$ws.Range("J10").Value = 0.5206762565675317
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 10.50918
$ws.Range("N10").Value = 31.52754
$ws.Range("O10").Value = 0.1224541750083835
$ws.Range("P10").Value = 0.1224541750083835
$ws.Range("Q10").Value = 3363.877337881719
$ws.Range("R10").Value = 30274.89604093548
$ws.Range("S10").Value = 0.06375898144443051
$ws.Range("T10").Value = 0.06375898144443051

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 320.0894206666666
$ws.Range("H11").Value = 960.2682619999999
$ws.Range("I11").Value = 0.5206762565675317
$ws.Range("J11").Value = 0.5206762565675317
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 5.176377666666667
$ws.Range("N11").Value = 15.529133
$ws.Range("O11").Value = 0.06031574839364136
$ws.Range("P11").Value = 0.06031574839364136
$ws.Range("Q11").Value = 1656.903728475205
$ws.Range("R11").Value = 14912.13355627685
$ws.Range("S11").Value = 0.0314049780856703
$ws.Range("T11").Value = 0.0314049780856703

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 320.0894206666666
$ws.Range("H12").Value = 960.2682619999999
$ws.Range("I12").Value = 0.5206762565675317
$ws.Range("J12").Value = 0.5206762565675317
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 66.063113
$ws.Range("N12").Value = 198.189339
$ws.Range("O12").Value = 0.7697749968028538
$ws.Range("P12").Value = 0.769774996802854
$ws.Range("Q12").Value = 21146.10356760653
$ws.Range("R12").Value = 190314.9321084588
$ws.Range("S12").Value = 0.4008035637345936
$ws.Range("T12").Value = 0.4008035637345937

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 320.0894206666666
$ws.Range("H13").Value = 960.2682619999999
$ws.Range("I13").Value = 0.5206762565675317
$ws.Range("J13").Value = 0.5206762565675317
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 4.072658
$ws.Range("N13").Value = 12.217974
$ws.Range("O13").Value = 0.04745507979512132
$ws.Range("P13").Value = 0.04745507979512132
$ws.Range("Q13").Value = 1303.614739793465
$ws.Range("R13").Value = 11732.53265814118
$ws.Range("S13").Value = 0.02470873330283728
$ws.Range("T13").Value = 0.02470873330283728

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 39.86989333333333
$ws.Range("H14").Value = 119.60968
$ws.Range("I14").Value = 0.06485471080959287
$ws.Range("J14").Value = 0.06485471080959287
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 10.50918
$ws.Range("N14").Value = 31.52754
$ws.Range("O14").Value = 0.1224541750083835
$ws.Range("P14").Value = 0.1224541750083835
$ws.Range("Q14").Value = 418.9998856207999
$ws.Range("R14").Value = 3770.9989705872
$ws.Range("S14").Value = 0.007941730107595987
$ws.Range("T14").Value = 0.007941730107595987

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 39.86989333333333
$ws.Range("H15").Value = 119.60968
$ws.Range("I15").Value = 0.06485471080959287
$ws.Range("J15").Value = 0.06485471080959287
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 5.176377666666667
$ws.Range("N15").Value = 15.529133
$ws.Range("O15").Value = 0.06031574839364136
$ws.Range("P15").Value = 0.06031574839364136
$ws.Range("Q15").Value = 206.3816254230489
$ws.Range("R15").Value = 1857.43462880744
$ws.Range("S15").Value = 0.003911760419333776
$ws.Range("T15").Value = 0.003911760419333776

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 39.86989333333333
$ws.Range("H16").Value = 119.60968
$ws.Range("I16").Value = 0.06485471080959287
$ws.Range("J16").Value = 0.06485471080959287
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 66.063113
$ws.Range("N16").Value = 198.189339
$ws.Range("O16").Value = 0.7697749968028538
$ws.Range("P16").Value = 0.769774996802854
$ws.Range("Q16").Value = 2633.929268577946
$ws.Range("R16").Value = 23705.36341720152
$ws.Range("S16").Value = 0.04992353480610436
$ws.Range("T16").Value = 0.04992353480610438

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 39.86989333333333
$ws.Range("H17").Value = 119.60968
$ws.Range("I17").Value = 0.06485471080959287
$ws.Range("J17").Value = 0.06485471080959287
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 4.072658
$ws.Range("N17").Value = 12.217974
$ws.Range("O17").Value = 0.04745507979512132
$ws.Range("P17").Value = 0.04745507979512132
$ws.Range("Q17").Value = 162.3764400431466
$ws.Range("R17").Value = 1461.38796038832
$ws.Range("S17").Value = 0.003077685476558747
$ws.Range("T17").Value = 0.003077685476558747
